$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet is a flat, formula-free price list. Rows 4..103 are 100 product
# rows (columns: A=#, B=name, H=counter code, L=price, N=counter code),
# row 104 is a totals row (K104 = sum of L), row 105 is a footer row
# (date / page / credit). A new product row ("OTAL EAR DROPS 5 ML") was
# inserted in alphabetical order right after row 59 ("ORS 10 SACHET"),
# pushing rows 60..103 down to 61..104, and the totals/footer rows down to
# 105/106.
# ---------------------------------------------------------------------------

# 1) Push the footer row (105 -> 106) down first, while row 105 still holds
#    the original footer content/formatting/merges.
$ws.Range("A105:N105").Copy($ws.Range("A106:N106"))

# 2) Push the totals row (104 -> 105) down next, while row 104 still holds
#    the original totals content/formatting/merges (blank A..J, K has the
#    sum, merged K:N).
$ws.Range("A104:N104").Copy($ws.Range("A105:N105"))

# 3) Turn (old) row 104 into a regular item row by cloning the formatting
#    and merge layout of the last item row (103).
$ws.Range("A103:N103").Copy($ws.Range("A104:N104"))

# 4) Capture the current (pre-shift) item values for rows 60..103 so we can
#    shift them down by one row into 61..104 without clobbering data we
#    still need to read.
$names = @{}
$codes1 = @{}
$prices = @{}
$codes2 = @{}
for ($r = 60; $r -le 103; $r++) {
    $names[$r]  = $ws.Range("B$r").Value2
    $codes1[$r] = $ws.Range("H$r").Value2
    $prices[$r] = $ws.Range("L$r").Value2
    $codes2[$r] = $ws.Range("N$r").Value2
}

for ($r = 103; $r -ge 61; $r--) {
    $src = $r - 1
    $ws.Range("B$r").Value = $names[$src]
    $ws.Range("H$r").Value = $codes1[$src]
    $ws.Range("L$r").Value = $prices[$src]
    $ws.Range("N$r").Value = $codes2[$src]
}

# 5) Write the brand-new product row at 60 (A60 keeps its original running
#    index of 57 - it's untouched).
$ws.Range("B60").Value = "OTAL EAR DROPS 5 ML"
$ws.Range("H60").Value = "4:0"
$ws.Range("L60").Value = 19
$ws.Range("N60").Value = "1:0"

# A104 is the new last item row's running index (100 -> 101).
$ws.Range("A104").Value = 101

# 6) Update the grand total (old 6107.35 + new row's price 19).
$ws.Range("K105").Value = 6126.3500000000004

Write-Output "done"
